$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '56.249.61'
$ws.Range('E2').Value = '  +2.84%  '
$ws.Range('D3').Value = '2.317.58'
$ws.Range('E3').Value = '  +1.52%  '
$ws.Range('E4').Value = '  +0.12%  '
$ws.Range('D5').Value = '''515.74'
$ws.Range('E5').Value = '  +1.82%  '
$ws.Range('D6').Value = '''133.11'
$ws.Range('E6').Value = '  +2.85%  '
$ws.Range('D7').Value = '''0.996'
$ws.Range('E7').Value = '  -0.01%  '
$ws.Range('D8').Value = '''0.534'
$ws.Range('E8').Value = '  +1.08%  '
$ws.Range('D9').Value = '2.337.33'
$ws.Range('D10').Value = '''0.102'
$ws.Range('E10').Value = '  +5.32%  '
$ws.Range('E11').Value = '  +0.17%  '
$ws.Range('E12').Value = '  +5.36%  '
$ws.Range('D13').Value = '''0.338'
$ws.Range('E13').Value = '  -2.09%  '
$ws.Range('D14').Value = '''23.61'
$ws.Range('E14').Value = '  +1.35%  '
$ws.Range('D15').Value = '2.734.93'
$ws.Range('E15').Value = '  +1.63%  '
$ws.Range('D16').Value = '56.483.02'
$ws.Range('E16').Value = '  +3.15%  '
$ws.Range('D17').Value = '''0.0000133'
$ws.Range('E17').Value = '  +1.68%  '
$ws.Range('D18').Value = '2.328.05'
$ws.Range('E18').Value = '  +1.10%  '
$ws.Range('D19').Value = '''10.36'
$ws.Range('E19').Value = '  +0.10%  '
$ws.Range('D20').Value = '''4.24'
$ws.Range('E20').Value = '  +2.20%  '
$ws.Range('D21').Value = '''318.57'
$ws.Range('E21').Value = '  +3.78%  '
$ws.Range('D22').Value = '''6.63'
$ws.Range('E22').Value = '  +3.39%  '
$ws.Range('D23').Value = '''0.998'
$ws.Range('E23').Value = '  +0.08%  '
$ws.Range('D24').Value = '''60.44'
$ws.Range('E24').Value = '  +0.35%  '
$ws.Range('E25').Value = '  +0.15%  '
$ws.Range('D26').Value = '''0.158'
$ws.Range('E26').Value = '  +4.36%  '
$ws.Range('D27').Value = '''7.72'
$ws.Range('E27').Value = '  +3.89%  '
$ws.Range('D28').Value = '''170.64'
$ws.Range('E28').Value = '  +0.03%  '
$ws.Range('D29').Value = '''1.20'
$ws.Range('E29').Value = '  +8.56%  '
$ws.Range('D30').Value = '0.0₃0731'
$ws.Range('E30').Value = '  +3.46%  '
$ws.Range('B31').Value = 'PancakeSwap'
$ws.Range('C31').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D31').Value = '''1.67'
$ws.Range('E31').Value = '  +2.01%  '
$ws.Range('B32').Value = 'Aptos'
$ws.Range('C32').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D32').Value = '''6.20'
$ws.Range('E32').Value = '  +2.15%  '
$ws.Range('D33').Value = '''18.22'
$ws.Range('E33').Value = '  +1.15%  '
$ws.Range('E34').Value = '  +0.02%  '
$ws.Range('D35').Value = '''0.993'
$ws.Range('E35').Value = '  -0.20%  '
$ws.Range('D36').Value = '''0.938'
$ws.Range('E36').Value = '  +1.52%  '
$ws.Range('D37').Value = '''1.24'
$ws.Range('E37').Value = '  +3.15%  '
$ws.Range('D38').Value = '''3.95'
$ws.Range('E38').Value = '  +4.35%  '
$ws.Range('E39').Value = '  +7.01%  '
$ws.Range('D40').Value = '''37.40'
$ws.Range('E40').Value = '  +2.73%  '
$ws.Range('D41').Value = '''0.378'
$ws.Range('E41').Value = '  +0.12%  '
$ws.Range('D42').Value = '''137.69'
$ws.Range('E42').Value = '  +9.76%  '
$ws.Range('D43').Value = '''3.55'
$ws.Range('E43').Value = '  +3.64%  '
$ws.Range('D44').Value = '''274.76'
$ws.Range('E44').Value = '  +9.61%  '
$ws.Range('D45').Value = '''5.04'
$ws.Range('E45').Value = '  -0.10%  '
$ws.Range('D46').Value = '''0.0926'
$ws.Range('E46').Value = '  +2.36%  '
$ws.Range('D47').Value = '''0.0503'
$ws.Range('E47').Value = '  +1.01%  '
$ws.Range('D48').Value = '''0.556'
$ws.Range('E48').Value = '  +1.15%  '
$ws.Range('E49').Value = '  +3.61%  '
$ws.Range('D50').Value = '''0.378'
$ws.Range('E50').Value = '  +0.54%  '
$ws.Range('D51').Value = '''16.72'
$ws.Range('E51').Value = '  +1.94%  '
